# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp banner in A1
$ws.Range("A1").Value = "Datos actualizados a 14 de Mayo de 2020 a las 07:35"

# Row 60: Kazajistan - refreshed counters
$ws.Range("B60").Value = 5571
$ws.Range("C60").Value = 154
$ws.Range("E60").Value = 3131

# Row 68: Hungria - refreshed counters
$ws.Range("B68").Value = 3380
$ws.Range("C68").Value = 39
$ws.Range("D68").Value = 1169
$ws.Range("E68").Value = 1775
$ws.Range("F68").Value = 49
$ws.Range("G68").Value = 6
$ws.Range("H68").Value = 436

# Row 75: Uzbekistan - refreshed counters
$ws.Range("B75").Value = 2620
$ws.Range("C75").Value = 8
$ws.Range("E75").Value = 533

# Rows 97/98: Kirguistan overtakes Hong Kong in the ranking (re-sorted by
# Casos totales, descending), so the two rows swap contents. Kirguistan
# also gets refreshed counters while Hong Kong keeps its prior figures.
$ws.Range("A97").Value = "Kirguistan"
$ws.Range("B97").Value = 1082
$ws.Range("C97").Value = 38
$ws.Range("D97").Value = 735
$ws.Range("E97").Value = 335
$ws.Range("F97").Value = 13
$ws.Range("G97").Value = 0
$ws.Range("H97").Value = 12

$ws.Range("A98").Value = "Hong Kong"
$ws.Range("B98").Value = 1051
$ws.Range("C98").Value = 0
$ws.Range("D98").Value = 1008
$ws.Range("E98").Value = 39
$ws.Range("F98").Value = 1
$ws.Range("G98").Value = 0
$ws.Range("H98").Value = 4
